# Add 2022-Q4 data
# 1) Insert a "2022-Q4" summary row into the "总计" sheet (shifting the
#    existing "2022-Q1" row down).
# 2) Insert a brand new "2022-Q4" worksheet (holding the fund-holding
#    detail rows for that quarter) positioned between "总计" and "2022-Q1".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: "总计" (summary) sheet - insert a new row 2 for 2022-Q4
# ---------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item(1)

$wsTotal.Rows.Item(2).Insert()
$wsTotal.Range("B2:D2").ClearFormats()

# Copy the index-column style (bold + border) from the row below onto
# the freshly inserted row so A2 matches the look of A3 (and the header).
$wsTotal.Range("A3").Copy()
$wsTotal.Range("A2").PasteSpecial(-4122)

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q4"
$wsTotal.Range("C2").Value = 5
$wsTotal.Range("D2").Value = 2.57

$wsTotal.Range("A3").Value = 1

# ---------------------------------------------------------------------
# Step 2: create the new "2022-Q4" worksheet (fund holdings detail)
# ---------------------------------------------------------------------
$wsQ1 = $wb.Worksheets.Item("2022-Q1")

# Duplicate the existing "2022-Q1" sheet so the new sheet inherits the
# same look & feel (fonts, borders, column layout), then place it right
# after "总计" (i.e. before "2022-Q1").
$wsQ1.Copy($null, $wsTotal)
$wsQ4 = $wb.Worksheets.Item(2)
$wsQ4.Name = "2022-Q4"

# The Q4 detail table has one extra data row (5 funds instead of 4) -
# append row 6 by duplicating the format of the last existing data row.
$wsQ4.Range("A5:H5").Copy()
$wsQ4.Range("A6").PasteSpecial(-4122)

# --- header row (unchanged text, already copied) ---

# --- row 2 ---
$wsQ4.Range("A2").Value = 0
$wsQ4.Range("B2").Value = "'720001"
$wsQ4.Range("C2").Value = "财通价值动量混合"
$wsQ4.Range("D2").Value = "'36.34"
$wsQ4.Range("E2").Value = "'79.52"
$wsQ4.Range("F2").Value = "'4.05"
$wsQ4.Range("G2").Value = "'1.4718"
$wsQ4.Range("H2").Value = 9

# --- row 3 ---
$wsQ4.Range("A3").Value = 1
$wsQ4.Range("B3").Value = "'001480"
$wsQ4.Range("C3").Value = "财通成长优选混合"
$wsQ4.Range("D3").Value = "'20.59"
$wsQ4.Range("E3").Value = "'91.01"
$wsQ4.Range("F3").Value = "'4.57"
$wsQ4.Range("G3").Value = "'0.9410"
$wsQ4.Range("H3").Value = 10

# --- row 4 ---
$wsQ4.Range("A4").Value = 2
$wsQ4.Range("B4").Value = "'009062"
$wsQ4.Range("C4").Value = "财通智慧成长混合A"
$wsQ4.Range("D4").Value = "'2.05"
$wsQ4.Range("E4").Value = "'86.49"
$wsQ4.Range("F4").Value = "'3.93"
$wsQ4.Range("G4").Value = "'0.0806"
$wsQ4.Range("H4").Value = 10

# --- row 5 ---
$wsQ4.Range("A5").Value = 3
$wsQ4.Range("B5").Value = "'009063"
$wsQ4.Range("C5").Value = "财通智慧成长混合C"
$wsQ4.Range("D5").Value = "'1.74"
$wsQ4.Range("E5").Value = "'86.49"
$wsQ4.Range("F5").Value = "'3.93"
$wsQ4.Range("G5").Value = "'0.0684"
$wsQ4.Range("H5").Value = 10

# --- row 6 ---
$wsQ4.Range("A6").Value = 4
$wsQ4.Range("B6").Value = "'002020"
$wsQ4.Range("C6").Value = "国都创新驱动灵活配置混合"
$wsQ4.Range("D6").Value = "'0.12"
$wsQ4.Range("E6").Value = "'83.47"
$wsQ4.Range("F6").Value = "'3.40"
$wsQ4.Range("G6").Value = "'0.0041"
$wsQ4.Range("H6").Value = 5

# Restore the original active sheet / selection (the "2022-Q1" sheet was
# the active tab before this edit, and stays that way afterwards). Re-fetch
# the sheet by name since the earlier Copy() operation shifted indices.
$wsQ1Final = $wb.Worksheets.Item("2022-Q1")
$wsQ1Final.Activate()
[void]$wsQ1Final.Range("A1").Select()

